$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1779566666666667
$ws.Range("H2").Value = 0.5338700000000001
$ws.Range("I2").Value = 0.01192558037548992
$ws.Range("J2").Value = 0.01192558037548992
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03057133333333333
$ws.Range("N2").Value = 0.09171399999999999
$ws.Range("O2").Value = 0.1084248755136686
$ws.Range("P2").Value = 0.1084248755136687
$ws.Range("Q2").Value = 0.005440372575555555
$ws.Range("R2").Value = 0.04896335318
$ws.Range("S2").Value = 0.001293029567640745
$ws.Range("T2").Value = 0.001293029567640745

$ws.Range("G3").Value = 0.1779566666666667
$ws.Range("H3").Value = 0.5338700000000001
$ws.Range("I3").Value = 0.01192558037548992
$ws.Range("J3").Value = 0.01192558037548992
$ws.Range("O3").Value = 0.8915751244863314
$ws.Range("P3").Value = 0.8915751244863314
$ws.Range("Q3").Value = 0.04473605188222223
$ws.Range("R3").Value = 0.40262446694
$ws.Range("S3").Value = 0.01063255080784918
$ws.Range("T3").Value = 0.01063255080784918

$ws.Range("I4").Value = 0.540575811616083
$ws.Range("J4").Value = 0.540575811616083
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03057133333333333
$ws.Range("N4").Value = 0.09171399999999999
$ws.Range("O4").Value = 0.1084248755136686
$ws.Range("P4").Value = 0.1084248755136687
$ws.Range("Q4").Value = 0.2466071862271111
$ws.Range("R4").Value = 2.219464676044
$ws.Range("S4").Value = 0.05861186508017419
$ws.Range("T4").Value = 0.05861186508017421

$ws.Range("I5").Value = 0.540575811616083
$ws.Range("J5").Value = 0.540575811616083
$ws.Range("O5").Value = 0.8915751244863314
$ws.Range("P5").Value = 0.8915751244863314
$ws.Range("S5").Value = 0.4819639465359088
$ws.Range("T5").Value = 0.4819639465359088

$ws.Range("G6").Value = 6.677692666666666
$ws.Range("I6").Value = 0.4474986080084269
$ws.Range("J6").Value = 0.4474986080084269
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03057133333333333
$ws.Range("N6").Value = 0.09171399999999999
$ws.Range("O6").Value = 0.1084248755136686
$ws.Range("P6").Value = 0.1084248755136687
$ws.Range("Q6").Value = 0.2041459684102222
$ws.Range("R6").Value = 1.837313715692
$ws.Range("S6").Value = 0.04851998086585369
$ws.Range("T6").Value = 0.0485199808658537

$ws.Range("G7").Value = 6.677692666666666
$ws.Range("I7").Value = 0.4474986080084269
$ws.Range("J7").Value = 0.4474986080084269
$ws.Range("O7").Value = 0.8915751244863314
$ws.Range("P7").Value = 0.8915751244863314
$ws.Range("S7").Value = 0.3989786271425733
$ws.Range("T7").Value = 0.3989786271425732

